# Fruta / hortaliza, semanal
#
# The weekly refresh re-shuffles the data rows (2-15) of the sheet: each
# row keeps its "dimension" columns (market/region/category/etc.) but the
# observation-specific columns (date, variety, quality, volume, prices,
# commercialization unit, price per kg and its unit) get redistributed
# across the 14 data rows according to a fixed permutation.
#
# Rather than hard-coding literal values (and risking mistakes with
# accented strings), we snapshot the current contents of the columns that
# change, then write them back out in the new row order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 15

# Columns whose *values* move between rows.
$cols = @(4, 8, 9, 10, 11, 12, 13, 14, 16, 17)   # D H I J K L M N P Q

# Maps new row number -> row number that currently holds the data that
# should end up there.
$rowSource = @{
    2  = 3
    3  = 14
    4  = 5
    5  = 6
    6  = 4
    7  = 11
    8  = 15
    9  = 12
    10 = 2
    11 = 13
    12 = 9
    13 = 7
    14 = 8
    15 = 10
}

# Snapshot every relevant cell value before we start overwriting anything.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Now write the values back out according to the permutation.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $src = $rowSource[$r]
    $rowVals = $snapshot[$src]
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c]
    }
}
